# Refactor; handle empty rows; columns
#
# Adds a new "FILLER" header column (I) to the sheet, mirroring the
# formatting of the existing header cells, and updates the current
# selection to the newly added cell below the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, bold/black
# header style) onto the new header cell I1, then set its text.
$null = $ws.Range("G1").Copy()
$null = $ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "FILLER"

# Move/collapse the active selection onto I2 (just below the new header)
$null = $ws.Range("I2").Select()
